# Added Port for positiontracking
# Fills in the new "positiontracking" row (B3/C3) in the port-mapping sheet
# and leaves the selection on the newly-entered cell, matching the
# spreadsheet's existing row layout/styling for columns B and C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "positiontracking"
$ws.Range("C3").Value = "server for positions"

$ws.Range("C3").Select() | Out-Null
